# daily auto push: 2026-01-10 18:43 UTC
# Insert two new data rows (2026/01/10 and 2026/01/11) just above the existing
# "2026/12/29" block, which starts at row 597. This shifts all subsequent rows
# down by two (old row 638 becomes row 640).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertAt = 597

# Insert two blank rows before row 597; everything from row 597 downward
# shifts down to rows 599+.
$ws.Range("A" + $insertAt + ":A" + ($insertAt + 1)).EntireRow.Insert()

# Helper: write a date-like text string (e.g. "2026/01/10") into a cell
# without letting Excel auto-convert it into a real date value/format.
# We stage the text in a scratch cell far outside the used range, force it
# to be stored as text, then copy/paste it into the destination (a
# copy/paste of an already-text cell keeps it as text without re-parsing
# it), and finally clear the scratch cell again.
$scratch = $ws.Cells.Item(100000, 1)
function Set-DateTextCell($row, $col, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value2 = $text
    $scratch.ClearFormats()
    $scratch.Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial(-4104) | Out-Null
}

function Set-TextCell($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value2 = $text
}

function Set-NumberCell($row, $col, $num) {
    $ws.Cells.Item($row, $col).Value2 = $num
}

# New row 597: 2026/01/10 (Sat)
Set-DateTextCell $insertAt 1 "2026/01/10"
Set-TextCell $insertAt 2 "土"
Set-NumberCell $insertAt 3 23
Set-NumberCell $insertAt 4 23

# New row 598: 2026/01/11 (Sun)
Set-DateTextCell ($insertAt + 1) 1 "2026/01/11"
Set-TextCell ($insertAt + 1) 2 "日"
Set-NumberCell ($insertAt + 1) 3 2
Set-NumberCell ($insertAt + 1) 4 18

# Clean up the scratch cell so it doesn't linger in the used range.
$scratch.ClearContents() | Out-Null
$excel.CutCopyMode = $false
